# "Add last couple of days of walks"
#
# The walk log lives on the "Sheet1" tab. Column G holds the extra distance
# logged for that period and column F is the running "actual" total
# (F = previous F + this row's G). The user has walked more in the last
# couple of days, so the latest entry (G2) grows from 4.5 to 12 - F2
# recalculates automatically from its existing formula (=F1+G2).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G2").Value = 12

# Leave the sheet the way the user would after typing the new figure -
# selection moves on to the next empty row in the actuals column.
$ws.Activate()
$ws.Range("G3").Select()
